# Applies the diff: convert F134:F145 text->number, and append the two new
# export batches (rows 146-157 numeric F, rows 158-169 text F) that bring
# the sheet from A1:F145 to A1:F169.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rows 134-145 already exist as inline-string text in column F.
# Re-enter the same numbers so they are stored as real numeric values.
$existingF = @{
    134 = 45
    135 = 5
    136 = 15
    137 = 15
    138 = 15
    139 = 45
    140 = 35
    141 = 15
    142 = 50
    143 = 20
    144 = 45
    145 = 30
}
foreach ($r in $existingF.Keys) {
    $ws.Cells.Item($r, 6).Value = $existingF[$r]
}

# --- Step 2: append rows 146-157 (13:15:05 batch) - column F numeric
$batch1 = @(
    @("2024-07-25", "13:15:05", "Epcot", "World Showcase", "Frozen Ever After", 45),
    @("2024-07-25", "13:15:05", "Epcot", "World Showcase", "Gran Fiesta Tour", 5),
    @("2024-07-25", "13:15:05", "Epcot", "World Discovery", "Guardians of the Galaxy: Cosmic Rewind", 15),
    @("2024-07-25", "13:15:05", "Epcot", "World Celebration", "Journey Into Imagination With Figment", 15),
    @("2024-07-25", "13:15:05", "Epcot", "World Nature", "Living with the Land", 15),
    @("2024-07-25", "13:15:05", "Epcot", "World Showcase", "Meet Anna and Elsa at Royal Sommerhus", 20),
    @("2024-07-25", "13:15:05", "Epcot", "World Celebration", "Meet Beloved Disney Pals at Mickey & Friends", 35),
    @("2024-07-25", "13:15:05", "Epcot", "World Discovery", "Mission: SPACE", 15),
    @("2024-07-25", "13:15:05", "Epcot", "World Showcase", "Remy's Ratatouille Adventure", 50),
    @("2024-07-25", "13:15:05", "Epcot", "World Nature", "The Seas with Nemo & Friends", 20),
    @("2024-07-25", "13:15:05", "Epcot", "World Nature", "Soarin'", 45),
    @("2024-07-25", "13:15:05", "Epcot", "World Celebration", "Spaceship Earth", 30),
)

$r = 146
foreach ($row in $batch1) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Step 3: append rows 158-169 (13:20:19 batch) - column F stored as text
$batch2 = @(
    @("2024-07-25", "13:20:19", "Epcot", "World Showcase", "Frozen Ever After", "45"),
    @("2024-07-25", "13:20:19", "Epcot", "World Showcase", "Gran Fiesta Tour", "5"),
    @("2024-07-25", "13:20:19", "Epcot", "World Discovery", "Guardians of the Galaxy: Cosmic Rewind", "15"),
    @("2024-07-25", "13:20:19", "Epcot", "World Celebration", "Journey Into Imagination With Figment", "15"),
    @("2024-07-25", "13:20:19", "Epcot", "World Nature", "Living with the Land", "10"),
    @("2024-07-25", "13:20:19", "Epcot", "World Showcase", "Meet Anna and Elsa at Royal Sommerhus", "20"),
    @("2024-07-25", "13:20:19", "Epcot", "World Celebration", "Meet Beloved Disney Pals at Mickey & Friends", "35"),
    @("2024-07-25", "13:20:19", "Epcot", "World Discovery", "Mission: SPACE", "15"),
    @("2024-07-25", "13:20:19", "Epcot", "World Showcase", "Remy's Ratatouille Adventure", "45"),
    @("2024-07-25", "13:20:19", "Epcot", "World Nature", "The Seas with Nemo & Friends", "20"),
    @("2024-07-25", "13:20:19", "Epcot", "World Nature", "Soarin'", "50"),
    @("2024-07-25", "13:20:19", "Epcot", "World Celebration", "Spaceship Earth", "30"),
)

$r = 158
foreach ($row in $batch2) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 6).Style = "Normal"
    $r = $r + 1
}
